# Update cryptos list with the latest scraped prices / volumes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.027.75'
$ws.Range('D3').Value = '1.918.57'
$ws.Range('E3').Value = '  +2.30%  '
$ws.Range('E4').Value = '  -0.68%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '327.01'
$ws.Range('E5').Value = '  +3.61%  '
$ws.Range('E6').Value = '  -0.60%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.5236'
$ws.Range('E7').Value = '  +3.17%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.4055'
$ws.Range('E8').Value = '  +3.72%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.08458'
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '42.93'
$ws.Range('E10').Value = '  +2.96%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '1.126'
$ws.Range('E11').Value = '  +2.07%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '22.04'
$ws.Range('E12').Value = '  +8.01%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '6.346'
$ws.Range('E13').Value = '  +1.92%  '
$ws.Range('D14').Value = '1.923.54'
$ws.Range('E14').Value = '  +2.53%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '7.368'
$ws.Range('E15').Value = '  +1.56%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '1.002'
$ws.Range('E16').Value = '  -0.80%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '96.05'
$ws.Range('E17').Value = '  +5.13%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.00001114'
$ws.Range('E18').Value = '  +0.98%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.06740'
$ws.Range('E19').Value = '  +0.21%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '18.20'
$ws.Range('E20').Value = '  +2.56%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '1.001'
$ws.Range('E21').Value = '  -0.58%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.056'
$ws.Range('E22').Value = '  +2.17%  '
$ws.Range('D23').Value = '30.038.24'
$ws.Range('E23').Value = '  +5.31%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '11.26'
$ws.Range('E24').Value = '  +1.35%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.199'
$ws.Range('E25').Value = '  -1.47%  '
$ws.Range('D26').Value = '2.144.80'
$ws.Range('E26').Value = '  +2.60%  '
$ws.Range('E27').Value = '  +2.72%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '160.19'
$ws.Range('E28').Value = '  -0.90%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.450'
$ws.Range('E29').Value = '  +2.57%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '129.28'
$ws.Range('E30').Value = '  +2.62%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.080'
$ws.Range('E31').Value = '  +3.78%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.1061'
$ws.Range('E32').Value = '  +1.45%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '6.075'
$ws.Range('E33').Value = '  +5.21%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '3.664'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.02517'
$ws.Range('E35').Value = '  +2.58%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.06606'
$ws.Range('E36').Value = '  +0.73%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.2220'
$ws.Range('E37').Value = '  +2.80%  '
$ws.Range('B38').Value = 'ARBITRUM'
$ws.Range('C38').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '1.236'
$ws.Range('E38').Value = '  +3.83%  '
$ws.Range('B39').Value = 'FraxShare'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '9.004'
$ws.Range('E39').Value = '  +2.09%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '5.193'
$ws.Range('E40').Value = '  +2.51%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.6561'
$ws.Range('E41').Value = '  +2.50%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '11.61'
$ws.Range('E42').Value = '  +4.62%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.242'
$ws.Range('E43').Value = '  -1.03%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.6175'
$ws.Range('E44').Value = '  +2.45%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '13.24'
$ws.Range('E45').Value = '  +1.23%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '3.754'
$ws.Range('E46').Value = '  +1.84%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.069'
$ws.Range('E47').Value = '  +3.13%  '
$ws.Range('B48').Value = 'Quant'
$ws.Range('C48').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '125.58'
$ws.Range('E48').Value = '  +3.04%  '
$ws.Range('B49').Value = 'EOS'
$ws.Range('C49').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.242'
$ws.Range('E49').Value = '  +2.05%  '
$ws.Range('E50').Value = '  +3.29%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '79.48'
$ws.Range('E51').Value = '  +4.06%  '
